# GP05MOAS-PG514 Cal Info update
# - Moorings sheet: Water Depth cell I2 becomes a plain number (1000) instead of text "1000 m"
# - Asset_Cal_Info sheet: FLORT calibration rows (6-9) were re-keyed under the correct
#   Ref Des/columns (shifted from A:G to G:M) and a note about the Summer 2015 glider
#   configuration (no FLORT) was added.

$wb = $excel.ActiveWorkbook

$moorings = $wb.Worksheets.Item("Moorings")
$moorings.Range("I2").Value = 1000

$cal = $wb.Worksheets.Item("Asset_Cal_Info")

# Note about the Summer 2015 profiling glider configuration (no FLORT installed)
$cal.Range("G5").Value = "Summer 2015 profiling glider has a special configuration and there is no flort"

# Move the FLORT (GP05MOAS-PG514-02-FLORTM000) calibration coefficient rows from
# columns A:G to columns G:M, so they line up under the correct headers.
$src = $cal.Range("A6:G9")
$dst = $cal.Range("G6:M9")
$src.Copy($dst)
$cal.Range("A6:F9").Clear()
